# "fixed test result write to xlsx file"
# Corrects the typo in B2, and writes the test result columns (Result / Pass)
# into column D for rows 1-2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo: "Hello  only me post" -> "Helo only me post"
$ws.Range("B2").Value = "Helo only me post"

# New "Result" column with the test outcome
$ws.Range("D1").Value = "Result"
$ws.Range("D2").Value = "Pass"

# Reflect the cursor position / zoom level that resulted from the edit
$ws.Range("D2").Select()
$excel.ActiveWindow.Zoom = 141
